$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.710.95'
$ws.Range("D3").Value = '3.548.90'
$ws.Range("E3").Value = '  -1.82%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '199.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '588.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.629'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.17'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.51%  '
$ws.Range("E12").Value = '  -4.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '686.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +15.68%  '
$ws.Range("D15").Value = '4.112.97'
$ws.Range("E15").Value = '  -1.92%  '
$ws.Range("D16").Value = '69.776.68'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").Value = '3.541.31'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("E18").Value = '  -5.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.17%  '
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.971'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '108.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.64%  '
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("E32").Value = '  -3.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.48%  '
$ws.Range("E34").Value = '  -4.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '62.43'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("D36").Value = '3.794.41'
$ws.Range("E36").Value = '  -3.65%  '
$ws.Range("D37").Value = '0.0₃0820'
$ws.Range("E37").Value = '  -6.72%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.74'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.07%  '
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '503.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.98%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.137'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("E43").Value = '  -5.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '34.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.26%  '
$ws.Range("E45").Value = '  +1.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("E50").Value = '  -2.57%  '
$ws.Range("E51").Value = '  +20.98%  '
